$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "Best Case Average Latency" values for the los/nlos test cases
$ws.Range("E7").Value = 8
$ws.Range("E13").Value = 5
$ws.Range("E19").Value = 7

# Update the active cell selection to match the final state
$ws.Range("E22").Select()
